$wb = $excel.ActiveWorkbook

# Update the daily conversion summary text on "Hoja1"
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 13057.3 pesos`n✅ 13057.3 pesos = 3.36 = 975.18 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the rate figures on "tasas"
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 297.55
$tasas.Range("O10").Value = 3885.2
$tasas.Range("N12").Value = 3883
$tasas.Range("O12").Value = 290
